# Applies the periodic cryptos.xlsx price/volume data refresh described by the commit.
# All Price (D) column values are plain text in this sheet (e.g. "1.861.23"), so we
# force a text number format before assigning, otherwise Excel COM auto-converts
# numeric-looking strings (e.g. "0.9997", "2.740", "1.000") into Double values and
# mangles/ truncates them (trailing zeros, exponential notation, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.296.93"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.61"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7058"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.44"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3143"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07813"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07988"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.86"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.01"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.194"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6984"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.425"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.311.98"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008260"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.23"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.111.56"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.586"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.0000"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1558"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.993"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.94"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.78"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.493"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.316"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.251"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05267"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.889"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7493"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.157"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.242.41"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.740"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8989"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.124"
$ws.Range("E42").Value = "  -6.60%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.98"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.80"
$ws.Range("E44").Value = "  -5.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000130"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.009.70"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5182"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.779"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.491"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4301"
$ws.Range("E51").Value = "  -1.77%  "
